$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with the day's new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.66 = 18384.54 pesos`n✅ 18384.54 pesos = 4.61 = 913.7 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: refresh the rate lookup cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 214.8
$wsTasas.Range("O10").Value = 3949
$wsTasas.Range("N12").Value = 3983.97
$wsTasas.Range("O12").Value = 198
